$d = $word.ActiveDocument

$replacements = @(
    @{old="42×79="; new="15×60="},
    @{old="77×60="; new="25×38="},
    @{old="58×71="; new="73×45="},
    @{old="45×42="; new="66×54="},
    @{old="45×26="; new="53×17="},
    @{old="29×14="; new="49×69="},
    @{old="47×86="; new="84×76="},
    @{old="46×35="; new="75×14="},
    @{old="76×45="; new="78×43="},
    @{old="52×16="; new="64×91="},
    @{old="99×31="; new="25×67="},
    @{old="89×16="; new="33×76="},
    @{old="27×23="; new="45×95="},
    @{old="55×49="; new="32×28="},
    @{old="72×75="; new="90×74="},
    @{old="48×16="; new="86×29="},
    @{old="60×19="; new="52×60="},
    @{old="44×58="; new="43×51="},
    @{old="50×20="; new="94×64="},
    @{old="78×15="; new="82×13="},
    @{old="83×34="; new="30×50="},
    @{old="59×66="; new="62×39="},
    @{old="72×47="; new="58×36="},
    @{old="66×53="; new="81×88="},
    @{old="28×89="; new="15×93="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
